# Update "Building Age Midpoint" (column C) summary-statistics values
# on the active worksheet to reflect corrected calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 11.2669322709163
$ws.Range("C3").Value  = 16.9295081967213
$ws.Range("C4").Value  = 34.8201948627104
$ws.Range("C5").Value  = 32.2272727272727
$ws.Range("C7").Value  = 21.399
$ws.Range("C10").Value = 7.42276422764228
$ws.Range("C11").Value = 11.7734375
$ws.Range("C14").Value = 8.96759259259259
$ws.Range("C15").Value = 11.4809160305344
$ws.Range("C16").Value = 17.0461095100865
$ws.Range("C17").Value = 24.8848684210526
